# "Removed Enter/Exit statuses on the same line."
# Row 2: Exchange -> ByBit, From date -> 44501, TP % -> 1
# Row 3: From date -> 44501, TP % -> 1, Strategy -> ScalpEmaRsiAdx_X
# Rows 4 & 5: these were duplicate "Enter"/"Exit" status rows for the same
# test line; clear them out entirely (only the date-formatted D/E cells
# remain, now blank) so the test grid effectively shrinks back to two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 changes
$ws.Range("B2").Value = "ByBit"
$ws.Range("D2").Value = 44501
$ws.Range("H2").Value = 1

# Row 3 changes
$ws.Range("D3").Value = 44501
$ws.Range("H3").Value = 1
$ws.Range("J3").Value = "ScalpEmaRsiAdx_X"

# Row 4 - drop every cell (formatting included) except D4/E4, whose
# contents are cleared but whose date number format is kept
$ws.Range("A4:C4").Clear()
$ws.Range("F4:J4").Clear()
$ws.Range("D4:E4").ClearContents()

# Row 5 - same cleanup
$ws.Range("A5:C5").Clear()
$ws.Range("F5:J5").Clear()
$ws.Range("D5:E5").ClearContents()

# Selection moved to D4
$ws.Range("D4").Select()
